# Update countries & provincias Spain
# Applies the COVID data refresh described in the commit: updated timestamp,
# updated numeric stats for several countries, and re-ordering of a few
# country pairs whose ranking changed (label + full row of stats moves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 14:52"

# --- Row 4: Estados Unidos (numeric refresh only) -------------------------
$ws.Range("B4").Value = 1064819
$ws.Range("C4").Value = 625
$ws.Range("D4").Value = 147473
$ws.Range("E4").Value = 855666
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 61680

# --- Row 19: India (numeric refresh only) ---------------------------------
$ws.Range("B19").Value = 33610
$ws.Range("C19").Value = 548
$ws.Range("E19").Value = 24094

# --- Rows 21/22: Ecuador & Portugal swap rank, Portugal gets new stats ----
$ws.Range("A21").Value = "Portugal"
$ws.Range("B21").Value = 25045
$ws.Range("C21").Value = 540
$ws.Range("D21").Value = 1519
$ws.Range("E21").Value = 22537
$ws.Range("F21").Value = 172
$ws.Range("G21").Value = 16
$ws.Range("H21").Value = 989

$ws.Range("A22").Value = "Ecuador"
$ws.Range("B22").Value = 24675
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 1557
$ws.Range("E22").Value = 22235
$ws.Range("F22").Value = 146
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 883

# --- Rows 28/29: Israel & Pakistan swap rank, Pakistan gets new stats -----
$ws.Range("A28").Value = "Pakistan"
$ws.Range("B28").Value = 16029
$ws.Range("C28").Value = 504
$ws.Range("D28").Value = 4052
$ws.Range("E28").Value = 11619
$ws.Range("F28").Value = 111
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = 358

$ws.Range("A29").Value = "Israel"
$ws.Range("B29").Value = 15870
$ws.Range("C29").Value = 36
$ws.Range("D29").Value = 8412
$ws.Range("E29").Value = 7239
$ws.Range("F29").Value = 117
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 219

# --- Rows 34/35: Polonia & Catar swap rank, Catar gets new stats ---------
$ws.Range("A34").Value = "Catar"
$ws.Range("B34").Value = 13409
$ws.Range("C34").Value = 845
$ws.Range("D34").Value = 1372
$ws.Range("E34").Value = 12027
$ws.Range("F34").Value = 72
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 10

$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 12781
$ws.Range("C35").Value = 141
$ws.Range("D35").Value = 3236
$ws.Range("E35").Value = 8917
$ws.Range("F35").Value = 160
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 628

# --- Row 54: Finlandia (numeric refresh only) ------------------------------
$ws.Range("E54").Value = 1984
$ws.Range("F54").Value = 48
$ws.Range("G54").Value = 5
$ws.Range("H54").Value = 211

# --- Rows 62/63: Tailandia & Barein swap rank, Barein gets new stats ------
$ws.Range("A62").Value = "Barein"
$ws.Range("B62").Value = 3037
$ws.Range("C62").Value = 116
$ws.Range("D62").Value = 1495
$ws.Range("E62").Value = 1534
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 8

$ws.Range("A63").Value = "Tailandia"
$ws.Range("B63").Value = 2954
$ws.Range("C63").Value = 7
$ws.Range("D63").Value = 2684
$ws.Range("E63").Value = 216
$ws.Range("F63").Value = 61
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 54

# --- Row 82: Republica de Macedonia (numeric refresh only) ----------------
$ws.Range("B82").Value = 1465
$ws.Range("C82").Value = 23
$ws.Range("D82").Value = 738
$ws.Range("E82").Value = 650
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 77

# --- Rows 181/182: Timor Oriental & Antigua y Barbuda swap rank ----------
$ws.Range("A181").Value = "Antigua y Barbuda"
$ws.Range("B181").Value = 24
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 11
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 1
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 3

$ws.Range("A182").Value = "Timor Oriental"
$ws.Range("B182").Value = 24
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 16
$ws.Range("E182").Value = 8
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0
